$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")

# Add new "commodity" and "year" columns read directly from config
$ws.Range("D1").Value = "commodity"
$ws.Range("E1").Value = "year"
$ws.Range("D2").Value = "light"
$ws.Range("E2").Value = 700
$ws.Range("E3").Value = 710
$ws.Range("E4").Value = 720

# Make "config" the active sheet/tab, with E7 selected
$ws.Activate() | Out-Null
$ws.Range("E7").Select() | Out-Null
